$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("MAZA RIOFRIO CINTHIA NATELAHI", 112),
    @("PANTA NIMA FREDDY ROLAND JUNIOR", 110),
    @("PALACIOS PANTA LUIS MIGUEL", 109),
    @("VEGA ZAPATA JESUS GABRIEL", 106),
    @("ELIAS MACHADO JUANA MARGOT", 104),
    @("CRISANTO CARMEN ROSITA ABIGAIL", 104),
    @("PANTA VARONA CANDY ELIZABETH", 104),
    @("SALAZAR VEGA MARIA FERNANDA", 103),
    @("HIDALGO MOSCOL YESSICA JAZMIN", 98),
    @("TALLEDO ELIAS ANDREA ALESSANDRA", 90)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
